# South Dakota 2017 MCAS workbook cleanup
# - rename header columns to snake_case field names
# - title-case the lowercase connector words ("de", "del", "la", "el")
#   inside a handful of municipality / state names
# - drop the trailing footnote rows (188:192 and 476:480)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Delete the trailing metadata rows -------------------------------------
# Delete bottom-most block first so the 188:192 row numbers below it are not
# shifted before we get to them.
$ws.Rows("476:480").Delete() | Out-Null
$ws.Rows("188:192").Delete() | Out-Null

# --- Header row --------------------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the "de"/"del"/"la"/"el" connectors in these names ----------
$ws.Range("B5").Value = "Pabellón De Arteaga"
$ws.Range("A36").Value = "Ciudad De México"
$ws.Range("A49").Value = "Estado De México"
$ws.Range("B49").Value = "Ecatepec De Morelos"
$ws.Range("B50").Value = "Ixtapan De La Sal"
$ws.Range("B63").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B68").Value = "Purísima Del Rincón"
$ws.Range("B71").Value = "Valle De Santiago"
$ws.Range("B74").Value = "Coyuca De Catalán"
$ws.Range("B76").Value = "Zihuatanejo De Azueta"
$ws.Range("B79").Value = "Taxco De Alarcón"
$ws.Range("B80").Value = "Técpan De Galeana"
$ws.Range("B85").Value = "Tulancingo De Bravo"
$ws.Range("B89").Value = "Autlán De Navarro"
$ws.Range("B95").Value = "San Martín De Bolaños"
$ws.Range("B96").Value = "San Miguel El Alto"
$ws.Range("B98").Value = "Tamazula De Gordiano"
$ws.Range("B99").Value = "Tepatitlán De Morelos"
$ws.Range("B146").Value = "Amealco De Bonfil"
$ws.Range("B153").Value = "Santa María Del Río"
$ws.Range("B157").Value = "Villa De Arista"
$ws.Range("B158").Value = "Villa De Ramos"
$ws.Range("B173").Value = "Ignacio De La Llave"
$ws.Range("B183").Value = "Tlaltenango De Sánchez Román"
